# Update the player roster table (A2:C18) on Sheet1.
# The edit reorders the players and swaps "Klay Thompson" (SG,SF / Dallas
# Mavericks) for "T.J. McConnell" (PG / Indiana Pacers).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @("Ja Morant",        "PG",          "Memphis Grizzlies"),
    @("Tyler Herro",      "PG,SG",       "Miami Heat"),
    @("Mikal Bridges",    "SG,SF,PF",    "New York Knicks"),
    @("Josh Giddey",      "PG,SG,SF",    "Chicago Bulls"),
    @("Gradey Dick",      "SG,SF",       "Toronto Raptors"),
    @("Miles Bridges",    "SF,PF",       "Charlotte Hornets"),
    @("Scottie Barnes",   "PG,SG,SF,PF", "Toronto Raptors"),
    @("Kyle Kuzma",       "PF",          "Washington Wizards"),
    @("Brook Lopez",      "C",           "Milwaukee Bucks"),
    @("Nikola Vucevic",   "PF,C",        "Chicago Bulls"),
    @("T.J. McConnell",   "PG",          "Indiana Pacers"),
    @("De'Aaron Fox",     "PG",          "Sacramento Kings"),
    @("Shaedon Sharpe",   "SG,SF",       "Portland Trail Blazers"),
    @("DeMar DeRozan",    "SF,PF",       "Sacramento Kings"),
    @("Evan Mobley",      "PF,C",        "Cleveland Cavaliers"),
    @("Luka Doncic",      "PG,SG",       "Dallas Mavericks"),
    @("Nick Richards",    "C",           "Phoenix Suns")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
